$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each row: sheet row number, Coin, Link, Price, Volume(1h)
$data = @(
    ,@(2, "Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "61.303.88", "  +7.83%  ", 0)
    ,@(3, "Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "3.343.47", "  +3.74%  ", 0)
    ,@(4, "TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.00", "  +0.08%  ", 1)
    ,@(5, "BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "412.26", "  +5.05%  ", 1)
    ,@(6, "Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "111.89", "  +4.64%  ", 1)
    ,@(7, "LidoStakedEther", "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth", "3.337.12", "  +3.85%  ", 0)
    ,@(8, "XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.568", "  -0.81%  ", 1)
    ,@(9, "USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "1.00", "  -0.03%  ", 1)
    ,@(10, "Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.625", "  +2.39%  ", 1)
    ,@(11, "Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.109", "  +14.37%  ", 1)
    ,@(12, "Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "38.95", "  +0.45%  ", 1)
    ,@(13, "TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.143", "  +0.83%  ", 1)
    ,@(14, "WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "3.891.82", "  +4.11%  ", 0)
    ,@(15, "Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "8.25", "  +2.22%  ", 1)
    ,@(16, "Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "19.13", "  +1.09%  ", 1)
    ,@(17, "WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "3.417.91", "  +6.16%  ", 0)
    ,@(18, "WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "61.099.36", "  +7.73%  ", 0)
    ,@(19, "Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "1.01", "  -1.82%  ", 1)
    ,@(20, "Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "10.54", "  -3.20%  ", 1)
    ,@(21, "ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.0000115", "  +9.41%  ", 1)
    ,@(22, "ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "3.26", "  -1.49%  ", 1)
    ,@(23, "BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "298.23", "  +0.84%  ", 1)
    ,@(24, "InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "12.34", "  -4.33%  ", 1)
    ,@(25, "Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "74.03", "  +0.43%  ", 1)
    ,@(26, "PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "3.10", "  -0.70%  ", 1)
    ,@(27, "EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "28.76", "  +3.63%  ", 1)
    ,@(28, "LEO", "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo", "4.49", "  +2.25%  ", 1)
    ,@(29, "RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "7.45", "  +3.68%  ", 1)
    ,@(30, "Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "7.55", "  -0.93%  ", 1)
    ,@(31, "Kaspa", "https://coinranking.com/coin/V8GxkwWow+kaspa-kas", "0.169", "  +0.43%  ", 1)
    ,@(32, "Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.110", "  +2.18%  ", 1)
    ,@(33, "Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "0.999", "  -0.05%  ", 1)
    ,@(34, "Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "11.27", "  +0.57%  ", 1)
    ,@(35, "Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "2.43", "  +14.97%  ", 1)
    ,@(36, "InjectiveProtocol", "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj", "39.50", "  +6.19%  ", 1)
    ,@(37, "VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.0476", "  -1.07%  ", 1)
    ,@(38, "OKB", "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb", "51.99", "  +0.66%  ", 1)
    ,@(39, "Stacks", "https://coinranking.com/coin/mMPrMcB7+stacks-stx", "3.10", "  +3.78%  ", 1)
    ,@(40, "FirstDigitalUSD", "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd", "0.994", "  -0.58%  ", 1)
    ,@(41, "LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "3.32", "  -5.98%  ", 1)
    ,@(42, "Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "139.04", "  +3.67%  ", 1)
    ,@(43, "Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.120", "  -0.03%  ", 1)
    ,@(44, "ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "1.88", "  -0.13%  ", 1)
    ,@(45, "TheGraph", "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt", "0.279", "  -0.25%  ", 1)
    ,@(46, "NEARProtocol", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", "3.83", "  -3.06%  ", 1)
    ,@(47, "Celestia", "https://coinranking.com/coin/YQcD0lBl7+celestia-tia", "16.35", "  -2.94%  ", 1)
    ,@(48, "RocketPoolETH", "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth", "4.007.44", "  +12.86%  ", 0)
    ,@(49, "WEMIXToken", "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix", "2.20", "  +4.76%  ", 1)
    ,@(50, "EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "21.26", "  -2.55%  ", 1)
    ,@(51, "Maker", "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr", "2.143.08", "  +0.02%  ", 0)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    # Column D (Price) stores plain text such as "1.00" or "61.303.88".
    # When the text looks like a valid number, force the cell to Text
    # format first so Excel keeps the literal string (trailing zeros,
    # thousands-separator dots) instead of silently coercing it to a number.
    if ($row[5] -eq 1) {
        $ws.Cells.Item($r, 4).NumberFormat = "@"
    }
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
